# Add the two new expense entries received from Soby to the "Internal" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Internal")

# Copy the date-formatted style from an existing entry (row 9, column E)
# onto the two new date cells before filling in their values, so they pick
# up the same number format (and style index) as the rest of the table.
$ws.Range("E9").Copy()
$ws.Range("E10:E11").PasteSpecial(-4122)  # xlPasteFormats

# Row 10: Sr No. 9 - Paid Procurement Charges for 2nd components and PCB
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 10000
$ws.Range("C10").Value = "Sobby"
$ws.Range("D10").Value = "Parag"
$ws.Range("E10").Value = 44523
$ws.Range("F10").Value = "Paid Procurement Charges for 2nd components and PCB"

# Row 11: Sr No. 10 - Paid for 2nd milestone partial
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 31100
$ws.Range("C11").Value = "Sobby"
$ws.Range("D11").Value = "Parag"
$ws.Range("E11").Value = 44241
$ws.Range("F11").Value = "Paid for 2nd milestone partial"

# Leave the selection on A12 as in the saved state of the edited workbook.
$ws.Range("A12").Select() | Out-Null

# The author finished up by switching back to the Expenses_proto1 tab and
# selecting cell I5 there before saving.
$ws1 = $wb.Worksheets.Item("Expenses_proto1")
$ws1.Activate() | Out-Null
$ws1.Range("I5").Select() | Out-Null
